$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 235.5
$ws.Range("I12").Value = 450
$ws.Range("K12").Value = 450
$ws.Range("M12").Value = -280

$ws.Range("H17").Value = 1996.1818
$ws.Range("J17").Value = 1996.1818
$ws.Range("L17").Value = 5988.5454
$ws.Range("N17").Value = -6324.5454

$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2674

$ws.Range("H33").Value = 160.28572
$ws.Range("I33").Value = 164.92308
$ws.Range("K33").Value = 164.92308
$ws.Range("M33").Value = 64.07692

$ws.Range("H40").Value = 2314.2856
$ws.Range("I40").Value = 2362.5
$ws.Range("J40").Value = 2250
$ws.Range("K40").Value = 2362.5
$ws.Range("L40").Value = 2250
$ws.Range("M40").Value = -2187.5
$ws.Range("N40").Value = -2600

$ws.Range("H86").Value = 4231.0713
$ws.Range("J86").Value = 3953.0833
$ws.Range("L86").Value = 3953.0833
$ws.Range("N86").Value = -6199.0833

$ws.Range("H89").Value = 4231.0713
$ws.Range("J89").Value = 3953.0833
$ws.Range("L89").Value = 19765.4165
$ws.Range("N89").Value = -30997.4165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents() | Out-Null

$ws.Range("H45").Value = 1626.25
$ws.Range("I45").Value = 1626.25
$ws.Range("K45").Value = 1626.25
$ws.Range("M45").Value = -1249.25

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents() | Out-Null

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents() | Out-Null

$ws.Range("H102").Value = 1736.8889
$ws.Range("I102").Value = 1703.875
$ws.Range("J102").Value = 2001
$ws.Range("K102").Value = 1703.875
$ws.Range("L102").Value = 2001
$ws.Range("M102").Value = -81.875
$ws.Range("N102").Value = -5245

$ws.Range("H122").Value = 335665.1
$ws.Range("I122").Value = 527655.4399999999
$ws.Range("J122").Value = 4045.4546
$ws.Range("K122").Value = 1582966.32
$ws.Range("L122").Value = 12136.3638
$ws.Range("M122").Value = -1580516.32
$ws.Range("N122").Value = -17036.3638

$ws.Range("H132").Value = 26658.5
$ws.Range("I132").Value = 2211.3333
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 6633.999899999999
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -4103.999899999999
$ws.Range("N132").Value = -305060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5847.5
$ws.Range("I86").Value = 4888
$ws.Range("J86").Value = 6807
$ws.Range("K86").Value = 4888
$ws.Range("L86").Value = 6807
$ws.Range("M86").Value = -3765
$ws.Range("N86").Value = -9053

$ws.Range("H89").Value = 5847.5
$ws.Range("I89").Value = 4888
$ws.Range("J89").Value = 6807
$ws.Range("K89").Value = 24440
$ws.Range("L89").Value = 34035
$ws.Range("M89").Value = -18824
$ws.Range("N89").Value = -45267

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4264.6313
$ws.Range("I58").Value = 1540.3334
$ws.Range("J58").Value = 5522
$ws.Range("K58").Value = 1540.3334
$ws.Range("L58").Value = 5522
$ws.Range("M58").Value = -1337.3334
$ws.Range("N58").Value = -5928

$ws.Range("H107").Value = 1196.5714
$ws.Range("I107").Value = 702.5
$ws.Range("J107").Value = 1394.2
$ws.Range("K107").Value = 702.5
$ws.Range("L107").Value = 1394.2
$ws.Range("M107").Value = 1217.5
$ws.Range("N107").Value = -5234.2

$ws.Range("H134").Value = 3372.2666
$ws.Range("I134").Value = 2560.4
$ws.Range("J134").Value = 4996
$ws.Range("K134").Value = 7681.200000000001
$ws.Range("L134").Value = 14988
$ws.Range("M134").Value = -5146.200000000001
$ws.Range("N134").Value = -20058

$ws.Range("H136").Value = 4264.6313
$ws.Range("I136").Value = 1540.3334
$ws.Range("J136").Value = 5522
$ws.Range("K136").Value = 4621.0002
$ws.Range("L136").Value = 16566
$ws.Range("M136").Value = -2071.0002
$ws.Range("N136").Value = -21666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1722.5385
$ws.Range("I117").Value = 303.75
$ws.Range("K117").Value = 911.25
$ws.Range("M117").Value = 2530.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3558.7646
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3558.7646
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3558.7646
$ws.Range("M102").ClearContents() | Out-Null
$ws.Range("N102").Value = -6802.7646

$ws.Range("H126").Value = 3641.4
$ws.Range("J126").Value = 3668.2222
$ws.Range("L126").Value = 11004.6666
$ws.Range("N126").Value = -15944.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4117.5454
$ws.Range("I16").Value = 4117.5454
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4117.5454
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3947.5454
$ws.Range("N16").ClearContents() | Out-Null

$ws.Range("H22").Value = 6566.3335
$ws.Range("I22").Value = 1099.5
$ws.Range("K22").Value = 1099.5
$ws.Range("M22").Value = -804.5

$ws.Range("H27").Value = 6566.3335
$ws.Range("I27").Value = 1099.5
$ws.Range("K27").Value = 1099.5
$ws.Range("M27").Value = -992.5

$ws.Range("H46").Value = 2884.077
$ws.Range("I46").Value = 1733.25
$ws.Range("J46").Value = 3870.5
$ws.Range("K46").Value = 1733.25
$ws.Range("L46").Value = 3870.5
$ws.Range("M46").Value = -1545.25
$ws.Range("N46").Value = -4246.5

$ws.Range("H127").Value = 52122.75
$ws.Range("J127").Value = 52122.75
$ws.Range("L127").Value = 52122.75
$ws.Range("N127").Value = -62042.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 555000
$ws.Range("J49").Value = 110000
$ws.Range("L49").Value = 110000
$ws.Range("N49").Value = -110460

$ws.Range("H50").Value = 24000
$ws.Range("J50").Value = 24000
$ws.Range("L50").Value = 24000
$ws.Range("N50").Value = -25262

$ws.Range("H64").Value = 60333.332
$ws.Range("J64").Value = 66500
$ws.Range("L64").Value = 66500
$ws.Range("N64").Value = -66996

$ws.Range("H67").Value = 60333.332
$ws.Range("J67").Value = 66500
$ws.Range("L67").Value = 66500
$ws.Range("N67").Value = -68216

$ws.Range("H136").Value = 1471.475
$ws.Range("I136").Value = 1185.5676
$ws.Range("J136").Value = 4997.6665
$ws.Range("K136").Value = 3556.7028
$ws.Range("L136").Value = 14992.9995
$ws.Range("M136").Value = -1006.7028
$ws.Range("N136").Value = -20092.9995

Write-Host "Done applying Seraph_Profits updates."
